$wb = $excel.ActiveWorkbook

# The "Swiss" sheet is the template for the new "Portugal" market sheet.
$swiss = $wb.Worksheets.Item("Swiss")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate Swiss (copies formatting, merges, column widths, styles) and
# place the copy after the last sheet.
$swiss.Copy($null, $lastSheet)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Update the market-specific values.
$portugal.Range("B4").Value = "NGC-3479/T2404"
$portugal.Range("B2").Value = "Portugal Market"

# Restore Swiss's selection to the full used range (deselecting it as the
# active tab), then make Portugal the active sheet/tab with B2 selected.
$swiss.Range("A1:D15").Select()
$portugal.Activate()
$portugal.Range("B2").Select()
